$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 5.423951666666667
$ws.Range("H2").Value = 16.271855
$ws.Range("I2").Value = 0.4774188439413272
$ws.Range("J2").Value = 0.4774188439413271
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 1.285895333333333
$ws.Range("N2").Value = 3.857686
$ws.Range("O2").Value = 0.3864259878905995
$ws.Range("P2").Value = 0.3864259878905995
$ws.Range("Q2").Value = 6.974634136392222
$ws.Range("R2").Value = 62.77170722752999
$ws.Range("S2").Value = 0.1844870484076153
$ws.Range("T2").Value = 0.1844870484076153
$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 5.423951666666667
$ws.Range("H3").Value = 16.271855
$ws.Range("I3").Value = 0.4774188439413272
$ws.Range("J3").Value = 0.4774188439413271
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 0.3517506666666667
$ws.Range("N3").Value = 1.055252
$ws.Range("O3").Value = 0.1057050253891921
$ws.Range("P3").Value = 0.1057050253891921
$ws.Range("Q3").Value = 1.907878614717778
$ws.Range("R3").Value = 17.17090753246
$ws.Range("S3").Value = 0.05046557102009674
$ws.Range("T3").Value = 0.05046557102009672
$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 5.423951666666667
$ws.Range("H4").Value = 16.271855
$ws.Range("I4").Value = 0.4774188439413272
$ws.Range("J4").Value = 0.4774188439413271
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 0.5966156666666667
$ws.Range("N4").Value = 1.789847
$ws.Range("O4").Value = 0.1792897076506553
$ws.Range("P4").Value = 0.1792897076506553
$ws.Range("Q4").Value = 3.236014539576111
$ws.Range("R4").Value = 29.124130856185
$ws.Range("S4").Value = 0.08559628495715436
$ws.Range("T4").Value = 0.08559628495715435
$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 5.423951666666667
$ws.Range("H5").Value = 16.271855
$ws.Range("I5").Value = 0.4774188439413272
$ws.Range("J5").Value = 0.4774188439413271
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 1.093401
$ws.Range("N5").Value = 3.280203
$ws.Range("O5").Value = 0.328579279069553
$ws.Range("P5").Value = 0.3285792790695531
$ws.Range("Q5").Value = 5.930554176284999
$ws.Range("R5").Value = 53.37498758656499
$ws.Range("S5").Value = 0.1568699395564608
$ws.Range("T5").Value = 0.1568699395564607
$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 1.583504333333333
$ws.Range("H6").Value = 4.750513
$ws.Range("I6").Value = 0.1393808158066948
$ws.Range("J6").Value = 0.1393808158066948
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 1.285895333333333
$ws.Range("N6").Value = 3.857686
$ws.Range("O6").Value = 0.3864259878905995
$ws.Range("P6").Value = 0.3864259878905995
$ws.Range("Q6").Value = 2.036220832546444
$ws.Range("R6").Value = 18.325987492918
$ws.Range("S6").Value = 0.05386036944109972
$ws.Range("T6").Value = 0.05386036944109971
$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 1.583504333333333
$ws.Range("H7").Value = 4.750513
$ws.Range("I7").Value = 0.1393808158066948
$ws.Range("J7").Value = 0.1393808158066948
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 0.3517506666666667
$ws.Range("N7").Value = 1.055252
$ws.Range("O7").Value = 0.1057050253891921
$ws.Range("P7").Value = 0.1057050253891921
$ws.Range("Q7").Value = 0.5569987049195556
$ws.Range("R7").Value = 5.012988344276
$ws.Range("S7").Value = 0.01473325267361298
$ws.Range("T7").Value = 0.01473325267361298
$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 1.583504333333333
$ws.Range("H8").Value = 4.750513
$ws.Range("I8").Value = 0.1393808158066948
$ws.Range("J8").Value = 0.1393808158066948
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 0.5966156666666667
$ws.Range("N8").Value = 1.789847
$ws.Range("O8").Value = 0.1792897076506553
$ws.Range("P8").Value = 0.1792897076506553
$ws.Range("Q8").Value = 0.9447434935012221
$ws.Range("R8").Value = 8.502691441510999
$ws.Range("S8").Value = 0.02498954571809214
$ws.Range("T8").Value = 0.02498954571809214
$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 1.583504333333333
$ws.Range("H9").Value = 4.750513
$ws.Range("I9").Value = 0.1393808158066948
$ws.Range("J9").Value = 0.1393808158066948
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 1.093401
$ws.Range("N9").Value = 3.280203
$ws.Range("O9").Value = 0.328579279069553
$ws.Range("P9").Value = 0.3285792790695531
$ws.Range("Q9").Value = 1.731405221571
$ws.Range("R9").Value = 15.582646994139
$ws.Range("S9").Value = 0.04579764797388994
$ws.Range("T9").Value = 0.04579764797388994
$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 0.6660723333333333
$ws.Range("H10").Value = 1.998217
$ws.Range("I10").Value = 0.05862800830537802
$ws.Range("J10").Value = 0.05862800830537802
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 1.285895333333333
$ws.Range("N10").Value = 3.857686
$ws.Range("O10").Value = 0.3864259878905995
$ws.Range("P10").Value = 0.3864259878905995
$ws.Range("Q10").Value = 0.8564993050957778
$ws.Range("R10").Value = 7.708493745862
$ws.Range("S10").Value = 0.02265538602746397
$ws.Range("T10").Value = 0.02265538602746397
$ws.Range("E11").Value = 3
$ws.Range("G11").Value = 0.6660723333333333
$ws.Range("H11").Value = 1.998217
$ws.Range("I11").Value = 0.05862800830537802
$ws.Range("J11").Value = 0.05862800830537802
$ws.Range("K11").Value = 3
$ws.Range("M11").Value = 0.3517506666666667
$ws.Range("N11").Value = 1.055252
$ws.Range("O11").Value = 0.1057050253891921
$ws.Range("P11").Value = 0.1057050253891921
$ws.Range("Q11").Value = 0.2342913872982222
$ws.Range("R11").Value = 2.108622485684
$ws.Range("S11").Value = 0.00619727510643775
$ws.Range("T11").Value = 0.006197275106437749
$ws.Range("E12").Value = 3
$ws.Range("G12").Value = 0.6660723333333333
$ws.Range("H12").Value = 1.998217
$ws.Range("I12").Value = 0.05862800830537802
$ws.Range("J12").Value = 0.05862800830537802
$ws.Range("K12").Value = 3
$ws.Range("M12").Value = 0.5966156666666667
$ws.Range("N12").Value = 1.789847
$ws.Range("O12").Value = 0.1792897076506553
$ws.Range("P12").Value = 0.1792897076506553
$ws.Range("Q12").Value = 0.3973891891998889
$ws.Range("R12").Value = 3.576502702799
$ws.Range("S12").Value = 0.01051139846921141
$ws.Range("T12").Value = 0.01051139846921141
$ws.Range("E13").Value = 3
$ws.Range("G13").Value = 0.6660723333333333
$ws.Range("H13").Value = 1.998217
$ws.Range("I13").Value = 0.05862800830537802
$ws.Range("J13").Value = 0.05862800830537802
$ws.Range("K13").Value = 3
$ws.Range("M13").Value = 1.093401
$ws.Range("N13").Value = 3.280203
$ws.Range("O13").Value = 0.328579279069553
$ws.Range("P13").Value = 0.3285792790695531
$ws.Range("Q13").Value = 0.7282841553389998
$ws.Range("R13").Value = 6.554557398050999
$ws.Range("S13").Value = 0.01926394870226488
$ws.Range("T13").Value = 0.01926394870226488
$ws.Range("E14").Value = 3
$ws.Range("G14").Value = 3.687463666666666
$ws.Range("H14").Value = 11.062391
$ws.Range("I14").Value = 0.3245723319466
$ws.Range("J14").Value = 0.3245723319466
$ws.Range("K14").Value = 3
$ws.Range("M14").Value = 1.285895333333333
$ws.Range("N14").Value = 3.857686
$ws.Range("O14").Value = 0.3864259878905995
$ws.Range("P14").Value = 0.3864259878905995
$ws.Range("Q14").Value = 4.741692320802889
$ws.Range("R14").Value = 42.675230887226
$ws.Range("S14").Value = 0.1254231840144205
$ws.Range("T14").Value = 0.1254231840144205
$ws.Range("E15").Value = 3
$ws.Range("G15").Value = 3.687463666666666
$ws.Range("H15").Value = 11.062391
$ws.Range("I15").Value = 0.3245723319466
$ws.Range("J15").Value = 0.3245723319466
$ws.Range("K15").Value = 3
$ws.Range("M15").Value = 0.3517506666666667
$ws.Range("N15").Value = 1.055252
$ws.Range("O15").Value = 0.1057050253891921
$ws.Range("P15").Value = 0.1057050253891921
$ws.Range("Q15").Value = 1.297067803059111
$ws.Range("R15").Value = 11.673610227532
$ws.Range("S15").Value = 0.03430892658904464
$ws.Range("T15").Value = 0.03430892658904463
$ws.Range("E16").Value = 3
$ws.Range("G16").Value = 3.687463666666666
$ws.Range("H16").Value = 11.062391
$ws.Range("I16").Value = 0.3245723319466
$ws.Range("J16").Value = 0.3245723319466
$ws.Range("K16").Value = 3
$ws.Range("M16").Value = 0.5966156666666667
$ws.Range("N16").Value = 1.789847
$ws.Range("O16").Value = 0.1792897076506553
$ws.Range("P16").Value = 0.1792897076506553
$ws.Range("Q16").Value = 2.199998593797444
$ws.Range("R16").Value = 19.799987344177
$ws.Range("S16").Value = 0.05819247850619735
$ws.Range("T16").Value = 0.05819247850619735
$ws.Range("E17").Value = 3
$ws.Range("G17").Value = 3.687463666666666
$ws.Range("H17").Value = 11.062391
$ws.Range("I17").Value = 0.3245723319466
$ws.Range("J17").Value = 0.3245723319466
$ws.Range("K17").Value = 3
$ws.Range("M17").Value = 1.093401
$ws.Range("N17").Value = 3.280203
$ws.Range("O17").Value = 0.328579279069553
$ws.Range("P17").Value = 0.3285792790695531
$ws.Range("Q17").Value = 4.031876460596999
$ws.Range("R17").Value = 36.286888145373
$ws.Range("S17").Value = 0.1066477428369375
$ws.Range("T17").Value = 0.1066477428369375
